$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Weekly update: add this week's new price reports for Alcachofa -------
# Two new observations arrived (Madrigal variety), inserted into the existing
# date-ordered list of records. All previously existing rows are preserved
# and simply shift down to make room.

# New record for 2023-05-23 (serial 45069) goes at the very top of the data
# block (row 8), pushing every existing record down by one row.
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 45069
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112013
$ws.Range("G8").Value = "Alcachofa"
$ws.Range("H8").Value = "Madrigal"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 17000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17500
$ws.Range("N8").Value = "`$/caja 40 unidades"
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 438
$ws.Range("Q8").Value = 40
$ws.Range("R8").Value = "Hortaliza"

# New record for 2023-05-24 (serial 45070) is inserted right after the old
# "44687" record (now at row 10), i.e. at row 11, pushing the remaining
# records down by one more row.
$ws.Rows("11:11").Insert()
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Terminal La Palmera de La Serena"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 45070
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 100112013
$ws.Range("G11").Value = "Alcachofa"
$ws.Range("H11").Value = "Madrigal"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 360
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("N11").Value = "`$/caja 40 unidades"
$ws.Range("O11").Value = "Provincia del Elquí"
$ws.Range("P11").Value = 438
$ws.Range("Q11").Value = 40
$ws.Range("R11").Value = "Hortaliza"

# Make sure both new date cells keep the same date/time number format used
# throughout the rest of column D.
$dateFormat = $ws.Range("D12").NumberFormat
$ws.Range("D8").NumberFormat = $dateFormat
$ws.Range("D11").NumberFormat = $dateFormat
